$wb = $excel.ActiveWorkbook

# Update "Elapsed Duration(Hrs)" column (G) values on several sheets to reflect
# the refreshed outage durations.

$wb.Worksheets.Item("R1").Range("G2").Value = "3926:00:22"
$wb.Worksheets.Item("R1").Range("G3").Value = "65:33:00"

$wb.Worksheets.Item("R2").Range("G2").Value = "12107:24:01"
$wb.Worksheets.Item("R2").Range("G3").Value = "3237:07:30"
$wb.Worksheets.Item("R2").Range("G4").Value = "475:19:04"

$wb.Worksheets.Item("R4").Range("G2").Value = "2953:13:50"
$wb.Worksheets.Item("R4").Range("G3").Value = "180:26:05"

$wb.Worksheets.Item("R5").Range("G2").Value = "427:12:49"

$wb.Worksheets.Item("R6").Range("G2").Value = "67:45:07"
